$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.444.53"
$ws.Range("E2").Value = "  -7.00%  "

$ws.Range("D3").Value = "2.205.55"
$ws.Range("E3").Value = "  -8.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.86"
$ws.Range("E5").Value = "  -3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.99"
$ws.Range("E6").Value = "  -14.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -11.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("E9").Value = "  -12.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.59"
$ws.Range("E10").Value = "  -13.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  -10.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.47"
$ws.Range("E12").Value = "  -14.22%  "

$ws.Range("E13").Value = "  -5.29%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.872"
$ws.Range("E14").Value = "  -13.87%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.535.72"
$ws.Range("E15").Value = "  -8.32%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.71"
$ws.Range("E16").Value = "  -13.96%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.195.37"
$ws.Range("E17").Value = "  -8.43%  "

$ws.Range("D18").Value = "42.311.36"
$ws.Range("E18").Value = "  -7.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.99"
$ws.Range("E19").Value = "  +4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.52"
$ws.Range("E20").Value = "  -13.30%  "

$ws.Range("D21").Value = "0.0₃0947"

$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.22"
$ws.Range("E22").Value = "  -10.31%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.80"
$ws.Range("E23").Value = "  -13.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "231.85"
$ws.Range("E24").Value = "  -12.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  -10.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  -11.18%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.55"
$ws.Range("E28").Value = "  -14.37%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("E29").Value = "  -8.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.38"
$ws.Range("E30").Value = "  -10.64%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0877"
$ws.Range("E31").Value = "  -10.11%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.02"
$ws.Range("E32").Value = "  -14.31%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.62"
$ws.Range("E33").Value = "  -10.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.68"
$ws.Range("E34").Value = "  -9.04%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.03"
$ws.Range("E35").Value = "  -2.04%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  -9.34%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  +5.17%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -12.75%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.24"
$ws.Range("E39").Value = "  -14.14%  "

$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  -14.41%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0314"
$ws.Range("E41").Value = "  -13.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").Value = "1.759.01"
$ws.Range("E43").Value = "  +4.90%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.99"
$ws.Range("E44").Value = "  -13.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.38"
$ws.Range("E45").Value = "  -15.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.204"
$ws.Range("E46").Value = "  -15.71%  "

$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.95"
$ws.Range("E47").Value = "  -16.25%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.23"
$ws.Range("E48").Value = "  -10.05%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.27"
$ws.Range("E49").Value = "  -16.10%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.50"
$ws.Range("E50").Value = "  -10.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.61"
$ws.Range("E51").Value = "  -13.63%  "

